# Update Sheets via scheduled runner
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 536.38464
$ws.Range("I39").Value = 196.625
$ws.Range("J39").Value = 1080
$ws.Range("K39").Value = 589.875
$ws.Range("L39").Value = 3240
$ws.Range("M39").Value = -293.875
$ws.Range("N39").Value = -3832
$ws.Range("H113").Value = 2681.2632
$ws.Range("I113").Value = 2475.5557
$ws.Range("J113").Value = 2866.4
$ws.Range("K113").Value = 2475.5557
$ws.Range("L113").Value = 2866.4
$ws.Range("M113").Value = 778.4443000000001
$ws.Range("N113").Value = -9374.4
$ws.Range("H125").Value = 1471.5
$ws.Range("I125").Value = 1043
$ws.Range("J125").Value = 1900
$ws.Range("K125").Value = 9387
$ws.Range("L125").Value = 17100
$ws.Range("M125").Value = -6927
$ws.Range("H129").Value = 418
$ws.Range("I129").Value = 324.33334
$ws.Range("J129").Value = 980
$ws.Range("K129").Value = 973.0000200000001
$ws.Range("L129").Value = 2940
$ws.Range("M129").Value = 4026.99998
$ws.Range("N129").Value = -12940
$ws.Range("H137").Value = 2518.9312
$ws.Range("I137").Value = 3088.303
$ws.Range("J137").Value = 1767.36
$ws.Range("K137").Value = 9264.909
$ws.Range("L137").Value = 5302.08
$ws.Range("M137").Value = -6714.909
$ws.Range("N137").Value = -10402.08

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 27804.973
$ws.Range("I32").Value = 4719.1455
$ws.Range("J32").Value = 107162.5
$ws.Range("K32").Value = 4719.1455
$ws.Range("L32").Value = 107162.5
$ws.Range("M32").Value = -4432.1455
$ws.Range("N32").Value = -107736.5
$ws.Range("H33").Value = 70000
$ws.Range("I33").Value = 70000
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 70000
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -69671
$ws.Range("N33").ClearContents()
$ws.Range("H61").Value = 1857
$ws.Range("I61").Value = 1608
$ws.Range("J61").Value = 2438
$ws.Range("K61").Value = 1608
$ws.Range("L61").Value = 2438
$ws.Range("M61").Value = -1396
$ws.Range("N61").Value = -2862
$ws.Range("H74").Value = 843.45
$ws.Range("I74").Value = 592.7857
$ws.Range("J74").Value = 1428.3334
$ws.Range("K74").Value = 592.7857
$ws.Range("L74").Value = 1428.3334
$ws.Range("M74").Value = 281.2143
$ws.Range("N74").Value = -3176.3334
$ws.Range("H77").Value = 843.45
$ws.Range("I77").Value = 592.7857
$ws.Range("J77").Value = 1428.3334
$ws.Range("K77").Value = 2963.9285
$ws.Range("L77").Value = 7141.666999999999
$ws.Range("M77").Value = 1404.0715
$ws.Range("N77").Value = -15877.667
$ws.Range("H110").Value = 753.1429000000001
$ws.Range("I110").Value = 766.4545000000001
$ws.Range("J110").Value = 704.3333
$ws.Range("K110").Value = 766.4545000000001
$ws.Range("L110").Value = 704.3333
$ws.Range("M110").Value = 1278.5455
$ws.Range("H122").Value = 3145.9443
$ws.Range("I122").Value = 2301.9285
$ws.Range("J122").Value = 6100
$ws.Range("K122").Value = 6905.7855
$ws.Range("L122").Value = 18300
$ws.Range("M122").Value = -4455.7855
$ws.Range("N122").Value = -23200
$ws.Range("H136").Value = 1857
$ws.Range("I136").Value = 1608
$ws.Range("J136").Value = 2438
$ws.Range("K136").Value = 4824
$ws.Range("L136").Value = 7314
$ws.Range("M136").Value = -2274
$ws.Range("N136").Value = -12414

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1255.7142
$ws.Range("I107").Value = 1058.5294
$ws.Range("J107").Value = 2093.75
$ws.Range("K107").Value = 1058.5294
$ws.Range("L107").Value = 2093.75
$ws.Range("M107").Value = 861.4706000000001
$ws.Range("N107").Value = -5933.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 31803.143
$ws.Range("I31").Value = 51010.43
$ws.Range("J31").Value = 2992.2144
$ws.Range("K31").Value = 51010.43
$ws.Range("L31").Value = 2992.2144
$ws.Range("M31").Value = -50715.43
$ws.Range("N31").Value = -3582.2144
$ws.Range("H34").Value = 31803.143
$ws.Range("I34").Value = 51010.43
$ws.Range("J34").Value = 2992.2144
$ws.Range("K34").Value = 51010.43
$ws.Range("L34").Value = 2992.2144
$ws.Range("M34").Value = -50808.43
$ws.Range("N34").Value = -3396.2144
$ws.Range("H132").Value = 2120.5
$ws.Range("I132").Value = 1692
$ws.Range("J132").Value = 3025.111
$ws.Range("K132").Value = 5076
$ws.Range("L132").Value = 9075.332999999999
$ws.Range("M132").Value = -2546
$ws.Range("N132").Value = -14135.333
$ws.Range("H110").Value = 30000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 30000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 30000
$ws.Range("N110").Value = -38180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 312.5
$ws.Range("I33").Value = 150
$ws.Range("J33").Value = 366.66666
$ws.Range("K33").Value = 900
$ws.Range("L33").Value = 2199.99996
$ws.Range("M33").Value = -617
$ws.Range("N33").Value = -2765.99996
$ws.Range("H36").Value = 2787.625
$ws.Range("I36").Value = 633.6667
$ws.Range("J36").Value = 4080
$ws.Range("K36").Value = 1901.0001
$ws.Range("L36").Value = 12240
$ws.Range("M36").Value = -1732.0001
$ws.Range("N36").Value = -12578
$ws.Range("H55").Value = 2237
$ws.Range("I55").Value = 1004
$ws.Range("J55").Value = 2909.5454
$ws.Range("K55").Value = 3012
$ws.Range("L55").Value = 8728.636200000001
$ws.Range("M55").Value = -2835
$ws.Range("N55").Value = -9082.636200000001
$ws.Range("H68").Value = 7985
$ws.Range("I68").Value = 318.36365
$ws.Range("J68").Value = 50151.5
$ws.Range("K68").Value = 955.09095
$ws.Range("L68").Value = 150454.5
$ws.Range("M68").Value = -144.09095
$ws.Range("N68").Value = -152076.5
$ws.Range("H71").Value = 7985
$ws.Range("I71").Value = 318.36365
$ws.Range("J71").Value = 50151.5
$ws.Range("K71").Value = 2865.27285
$ws.Range("L71").Value = 451363.5
$ws.Range("M71").Value = 1190.72715
$ws.Range("N71").Value = -459475.5
$ws.Range("H122").Value = 1207.3549
$ws.Range("I122").Value = 944.5
$ws.Range("J122").Value = 1685.2727
$ws.Range("K122").Value = 8500.5
$ws.Range("L122").Value = 15167.4543
$ws.Range("M122").Value = -6050.5
$ws.Range("N122").Value = -20067.4543
$ws.Range("H131").Value = 779.86487
$ws.Range("I131").Value = 440
$ws.Range("J131").Value = 905.7406999999999
$ws.Range("K131").Value = 1320
$ws.Range("L131").Value = 2717.2221
$ws.Range("M131").Value = 3720
$ws.Range("N131").Value = -12797.2221

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2337.7368
$ws.Range("I122").Value = 1219.909
$ws.Range("J122").Value = 3874.75
$ws.Range("K122").Value = 3659.727
$ws.Range("L122").Value = 11624.25
$ws.Range("M122").Value = -1209.727
$ws.Range("N122").Value = -16524.25
$ws.Range("H132").Value = 2884.4546
$ws.Range("I132").Value = 2524.7856
$ws.Range("J132").Value = 3513.875
$ws.Range("K132").Value = 7574.3568
$ws.Range("L132").Value = 10541.625
$ws.Range("M132").Value = -5044.3568
$ws.Range("N132").Value = -15601.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1514.3529
$ws.Range("I7").Value = 1119.5
$ws.Range("J7").Value = 2462
$ws.Range("K7").Value = 1119.5
$ws.Range("L7").Value = 2462
$ws.Range("M7").Value = -1007.5
$ws.Range("N7").Value = -2686
$ws.Range("H40").Value = 2410.0454
$ws.Range("I40").Value = 2417.75
$ws.Range("J40").Value = 2333
$ws.Range("K40").Value = 2417.75
$ws.Range("L40").Value = 2333
$ws.Range("M40").Value = -2281.75
$ws.Range("N40").Value = -2605
$ws.Range("H46").Value = 979500
$ws.Range("I46").Value = 7500
$ws.Range("J46").Value = 1951500
$ws.Range("K46").Value = 7500
$ws.Range("L46").Value = 1951500
$ws.Range("M46").Value = -7312
$ws.Range("N46").Value = -1951876
$ws.Range("H122").Value = 51110.523
$ws.Range("I122").Value = 86825.914
$ws.Range("J122").Value = 3490
$ws.Range("K122").Value = 260477.742
$ws.Range("L122").Value = 10470
$ws.Range("M122").Value = -258027.742
$ws.Range("N122").Value = -15370
$ws.Range("H126").Value = 1514.3529
$ws.Range("I126").Value = 1119.5
$ws.Range("J126").Value = 2462
$ws.Range("K126").Value = 3358.5
$ws.Range("L126").Value = 7386
$ws.Range("M126").Value = -888.5
$ws.Range("N126").Value = -12326
$ws.Range("H132").Value = 5086.9487
$ws.Range("I132").Value = 2871.1667
$ws.Range("J132").Value = 8632.200000000001
$ws.Range("K132").Value = 8613.500100000001
$ws.Range("L132").Value = 25896.6
$ws.Range("M132").Value = -6083.500100000001
$ws.Range("N132").Value = -30956.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 371169.75
$ws.Range("I126").Value = 500599.4
$ws.Range("J126").Value = 1370.7142
$ws.Range("K126").Value = 1501798.2
$ws.Range("L126").Value = 4112.142599999999
$ws.Range("M126").Value = -1499328.2
$ws.Range("N126").Value = -9052.142599999999
$ws.Range("H132").Value = 1150.1072
$ws.Range("I132").Value = 705.625
$ws.Range("J132").Value = 1742.75
$ws.Range("K132").Value = 2116.875
$ws.Range("L132").Value = 5228.25
$ws.Range("M132").Value = 413.125
$ws.Range("N132").Value = -10288.25

Write-Output "Applied diff changes"